# Command 0x0500 / ADXL345 bring-up update
# - Adds "ADXL345 Pin Out" and "ADXL345 Reg Spec" sheets between
#   "TIVA Pin OUT" and "Commands".
# - Fills in the ADXL345 pinout table (14-pin breakout) with two
#   formula-driven "Note" cells that cross-reference the I2C pins on
#   the "TIVA Pin OUT" sheet.
# - Marks the "ADXL345 Reg Spec" sheet as TBD (register map still to be
#   documented).
# - Re-selects the header row on "TIVA Pin OUT" (whole-row selection)
#   and leaves it as the active sheet/tab, matching the saved view
#   state from the author's Excel session.

$wb = $excel.ActiveWorkbook

$tiva = $wb.Worksheets.Item("TIVA Pin OUT")

# --- Create the two new sheets, in order, right after "TIVA Pin OUT" ---
$adxlPins = $wb.Worksheets.Add($null, $tiva)
$adxlPins.Name = "ADXL345 Pin Out"

$adxlReg = $wb.Worksheets.Add($null, $adxlPins)
$adxlReg.Name = "ADXL345 Reg Spec"

# --- "ADXL345 Pin Out" content -----------------------------------------
# Cells are written in this particular order so that new shared-string
# entries are introduced in the same sequence as the source edit.
$adxlPins.Range("B1").Value = 'I2C use'
$adxlPins.Range("C1").Value = 'Name on board'
$adxlPins.Range("D1").Value = 'QuadCopter Pin'
$adxlPins.Range("D8").Value = 'VDD'
$adxlPins.Range("E1").Value = 'Note'
$adxlPins.Range("B2").Value = 'vdd'
$adxlPins.Range("B3").Value = 'gnd'
$adxlPins.Range("B4").Value = 'rfu'
$adxlPins.Range("B7").Value = 'vs'
$adxlPins.Range("B8").Value = 'cs'
$adxlPins.Range("C8").Value = 'cs!'
$adxlPins.Range("D13").Value = 'GND'
$adxlPins.Range("E13").Value = 'means the address is 0x53 (If we tie it to VDD we would use a different address)'
$adxlPins.Range("B9").Value = 'INT1'
$adxlPins.Range("B10").Value = 'INT2'
$adxlPins.Range("B11").Value = 'NC'
$adxlPins.Range("B13").Value = 'Alt Address'
$adxlPins.Range("C13").Value = 'sdo'
$adxlPins.Range("C14").Value = 'sda'
$adxlPins.Range("B15").Value = 'scl/sclk'
$adxlPins.Range("C15").Value = 'scl'
$adxlPins.Range("E8").Value = 'This selects I2C as the comm (GND means it is set for SPI)'

# --- "ADXL345 Reg Spec" content (register map still TBD) ---------------
$adxlReg.Range("A1").Value = 'TBD'

# --- Remaining "ADXL345 Pin Out" cells (reuse strings already added) ---
$adxlPins.Range("A1").Value = 'Pin'
$adxlPins.Range("B5").Value = 'gnd'
$adxlPins.Range("B6").Value = 'gnd'
$adxlPins.Range("B12").Value = 'rfu'
$adxlPins.Range("B14").Value = 'Data'

$adxlPins.Range("A2").Value = 1
$adxlPins.Range("A3").Value = 2
$adxlPins.Range("A4").Value = 3
$adxlPins.Range("A5").Value = 4
$adxlPins.Range("A6").Value = 5
$adxlPins.Range("A7").Value = 6
$adxlPins.Range("A8").Value = 7
$adxlPins.Range("A9").Value = 8
$adxlPins.Range("A10").Value = 9
$adxlPins.Range("A11").Value = 10
$adxlPins.Range("A12").Value = 11
$adxlPins.Range("A13").Value = 12
$adxlPins.Range("A14").Value = 13
$adxlPins.Range("A15").Value = 14

# Pins 13/14 (SDA/SCL) echo the I2C pin assignments already documented
# on the "TIVA Pin OUT" sheet.
$adxlPins.Range("D14").Formula = "='TIVA Pin OUT'!E5"
$adxlPins.Range("D15").Formula = "='TIVA Pin OUT'!E4"

# --- Column widths (best-fit-ish) for "ADXL345 Pin Out" -----------------
$adxlPins.Columns.Item(1).ColumnWidth = 2.917466666666667
$adxlPins.Columns.Item(2).ColumnWidth = 16.917466666666666
$adxlPins.Columns.Item(3).ColumnWidth = 13.7508
$adxlPins.Columns.Item(4).ColumnWidth = 21.417466666666666
$adxlPins.Columns.Item(5).ColumnWidth = 9.084133333333334
$adxlPins.Columns.Item(6).ColumnWidth = 15.917466666666666
$adxlPins.Columns.Item(7).ColumnWidth = 11.2508
$adxlPins.Columns.Item(8).ColumnWidth = 14.417466666666668

# --- Restore "TIVA Pin OUT" as the active sheet/tab with a whole-row
#     selection on row 1 (matches the author's final saved view state).
$tiva.Activate()
$tiva.Range("A1:XFD1").Select()
